$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 49/50: coin order swap (MultiversX now ranked above HuobiToken)
$ws.Range("B49").Value = "MultiversX"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.00"
$ws.Range("E49").Value = "  -0.19%  "

$ws.Range("B50").Value = "HuobiToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.84"
$ws.Range("E50").Value = "  +8.64%  "

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.888.79"
$ws.Range("E2").Value = "  -0.34%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.308.23"
$ws.Range("E3").Value = "  +0.36%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.56"
$ws.Range("E5").Value = "  +2.47%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.50"
$ws.Range("E6").Value = "  -1.35%  "

# Row 7
$ws.Range("E7").Value = "  -2.15%  "

# Row 8
$ws.Range("E8").Value = "  +0.03%  "

# Row 9
$ws.Range("E9").Value = "  -2.42%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.46"
$ws.Range("E10").Value = "  -1.65%  "

# Row 11
$ws.Range("E11").Value = "  +0.22%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "18.61"
$ws.Range("E12").Value = "  +5.06%  "

# Row 13
$ws.Range("E13").Value = "  +1.28%  "

# Row 14
$ws.Range("E14").Value = "  -1.33%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.665.94"
$ws.Range("E15").Value = "  +0.31%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.302.29"
$ws.Range("E16").Value = "  +2.09%  "

# Row 17
$ws.Range("E17").Value = "  -0.50%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.810.30"
$ws.Range("E18").Value = "  -0.28%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.10"
$ws.Range("E19").Value = "  +2.72%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0899"
$ws.Range("E20").Value = "  -1.06%  "

# Row 21
$ws.Range("E21").Value = "  -1.23%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.42"
$ws.Range("E22").Value = "  -1.98%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.43"
$ws.Range("E23").Value = "  -0.60%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.14"
$ws.Range("E24").Value = "  -0.58%  "

# Row 25
$ws.Range("E25").Value = "  +0.89%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.02%  "

# Row 27
$ws.Range("E27").Value = "  +0.07%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "25.27"
$ws.Range("E28").Value = "  +1.24%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.39"
$ws.Range("E29").Value = "  +16.71%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "165.79"
$ws.Range("E30").Value = "  +0.66%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.08"
$ws.Range("E31").Value = "  -0.43%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "33.17"
$ws.Range("E32").Value = "  +0.46%  "

# Row 33
$ws.Range("E33").Value = "  +0.05%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.79"
$ws.Range("E34").Value = "  -0.26%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.84"
$ws.Range("E36").Value = "  -0.29%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.40"
$ws.Range("E37").Value = "  -0.57%  "

# Row 38
$ws.Range("E38").Value = "  -0.41%  "

# Row 39
$ws.Range("E39").Value = "  -1.08%  "

# Row 40
$ws.Range("E40").Value = "  -0.70%  "

# Row 41
$ws.Range("E41").Value = "  -0.96%  "

# Row 42
$ws.Range("E42").Value = "  -2.55%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.012.75"
$ws.Range("E43").Value = "  -0.13%  "

# Row 44
$ws.Range("E44").Value = "  -2.10%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.25"
$ws.Range("E45").Value = "  +4.15%  "

# Row 46
$ws.Range("E46").Value = "  -3.24%  "

# Row 47
$ws.Range("E47").Value = "  -6.56%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.81"
$ws.Range("E48").Value = "  -0.65%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.530.44"
$ws.Range("E51").Value = "  +0.11%  "

